$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A25").NumberFormat = "@"
$ws.Range("A25").Value = "2025-09-09"
$ws.Range("A25").Style = "Normal"
$ws.Range("B25").Value = 57.06000137329102
$ws.Range("C25").Value = 715.5499877929688
$ws.Range("D25").Value = 325.7000122070312
